# nevergiveup/Excel/Stage_关卡表.xlsx - "feat：some bgm sound fix"
#
# Lower the bgm volume (column V, "bgm大小，默认0.5") from 1 to 0.6 for the
# stage rows 21-80, and restore the sheet's last-used window
# position/selection (scrolled down toward W76) the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data fix: bgm volume 1 -> 0.6 for rows 21 through 80 (column V) ---
for ($row = 21; $row -le 80; $row++) {
    $ws.Cells.Item($row, 22).Value = 0.6
}

# --- restore the view/selection state recorded in the saved workbook ---
$win = $excel.ActiveWindow

# Keep the existing freeze (top 3 rows stay frozen) and just move the
# selection the way it ended up after scrolling the sheet.
[void]$ws.Range("W76").Select()
$win.ScrollRow = 60
$win.ScrollColumn = 17
